# Adding Example Data
# Replaces the sample row 2 with "Welcome" example data and fills in the
# previously-empty row 3 with a second "Users" example row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : "Welcome" example -------------------------------------------
# Write the new shared-string-backed cells in the same order the strings
# first appear so the rebuilt sharedStrings table lines up with the target:
#   ./asset/image/example.image.jpg, #FFFFFF, Welcome, Users, #EC3C66
$ws.Range("I2").Value = "./asset/image/example.image.jpg"
$ws.Range("H2").Value = "#FFFFFF"
$ws.Range("B2").Value = "Welcome"

# --- Row 3 : "Users" example (was completely empty before) ---------------
$ws.Range("B3").Value = "Users"

# Now the rest of the shared/colour string, reusing what already exists.
$ws.Range("E2").Value = "#EC3C66"
$ws.Range("E3").Value = "#EC3C66"
$ws.Range("H3").Value = "#FFFFFF"
$ws.Range("I3").Value = "./asset/image/example.image.jpg"
$ws.Range("J2").Value = "./asset/font/Sportage-DemoItalic.otf"
$ws.Range("J3").Value = "./asset/font/Sportage-DemoItalic.otf"

# --- Numeric cells ---------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 100
$ws.Range("F2").Value = 150
$ws.Range("G2").Value = 10

$ws.Range("A3").Value = 1
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = 300
$ws.Range("F3").Value = 120
$ws.Range("G3").Value = 5

# --- Selection / scroll position -------------------------------------------
$ws.Range("G2").Select()
try { $excel.ActiveWindow.ScrollColumn = 2 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("B1") } catch {}

# --- Workbook window size/position (best effort) ---------------------------
try {
    $win = $excel.ActiveWindow
    $win.Left = 2620
    $win.Top = 2620
    $win.Width = 14400
    $win.Height = 7360
} catch {}
